$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.362450838088989
$ws.Range("B1").Value = 2.270521879196167
$ws.Range("C1").Value = 2.716475486755371
$ws.Range("D1").Value = 3.131510257720947
$ws.Range("E1").Value = 2.308110952377319
